$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D12").Value = "ML & DL Book Landscape 마우스패드 이벤트 안내"
$ws.Range("E12").Value = "https://tensorflow.blog/2022/07/05/ml-dl-book-landscape-%eb%a7%88%ec%9a%b0%ec%8a%a4%ed%8c%a8%eb%93%9c-%ec%9d%b4%eb%b2%a4%ed%8a%b8-%ec%95%88%eb%82%b4/"

$ws.Range("D37").Value = "dsba_seminar"

$ws.Range("D45").Value = "Anomaly detection - AutoEncoder 활용"
$ws.Range("E45").Value = "https://dive-into-ds.tistory.com/112"

$ws.Range("D50").Value = "허준이 교수: 한국시간으로 오후 4시, 2022년 7월5일 필즈 메달 수상"
$ws.Range("E50").Value = "http://incredible.egloos.com/7544067"

$ws.Range("D51").Value = "[matplotlib] 그래프 x축, y축 눈금 위치 설정하는 방법, xticks, yticks"
$ws.Range("E51").Value = "https://bskyvision.com/entry/matplotlib-%EA%B7%B8%EB%9E%98%ED%94%84-x%EC%B6%95-y%EC%B6%95-%EB%88%88%EA%B8%88-%EC%9C%84%EC%B9%98-%EC%84%A4%EC%A0%95%ED%95%98%EB%8A%94-%EB%B0%A9%EB%B2%95-xticks-yticks"
